$d = $word.ActiveDocument

# =====================================================================
# PART 1: paragraph 2 ("L’application GlouGlou a pour but ...")
#   - merge "L’application " + "GlouGlou" + " a " into a single run
#     (drops the spell-check proofErr markers around GlouGlou)
#   - merge "pour but ... l’" + "apprentissage ... ce " into a single
#     run (drops the _GoBack bookmark that used to sit between them)
# =====================================================================

$p2Start = $d.Paragraphs(2).Range.Start

$cut1 = $p2Start + "L’application GlouGlou a ".Length
$cut2 = $cut1 + "pour but de présenter les problèmes de mélange de manière intuitive et compréhensible afin d’apporter une aide à l’apprentissage de la résolution de ce ".Length

$tb1 = $d.Bookmarks.Add("ZZTB1", $d.Range($cut1, $cut1))
$tb2 = $d.Bookmarks.Add("ZZTB2", $d.Range($cut2, $cut2))

$rA = $d.Range($p2Start, $cut1)
$rA.Text = "L’application GlouGlou a  "
$d.Range($cut1, $cut1 + 1).Delete()

$d.Bookmarks("_GoBack").Delete()

$rB = $d.Range($cut1, $cut2)
$rB.Text = $rB.Text + " "
$d.Range($cut2, $cut2 + 1).Delete()

$d.Bookmarks("ZZTB1").Delete()
$d.Bookmarks("ZZTB2").Delete()

Write-Output ("Para2: [" + $d.Paragraphs(2).Range.Text + "]")

# =====================================================================
# PART 2: paragraph 3 ("Objectifs") -> empty paragraph
# =====================================================================

$p3 = $d.Paragraphs(3)
$clr = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$clr.Text = ""

Write-Output ("Para3: [" + $d.Paragraphs(3).Range.Text + "]")

# =====================================================================
# PART 3: paragraph 4 ("Fournir une première approche ... utilisation")
#   split into 4 runs + reinsert _GoBack bookmark between "mélanges" and
#   " grâce"
# =====================================================================

$p4Start = $d.Paragraphs(4).Range.Start

$s1 = "Fournir une première approche ludique dans la rés"
$s2 = "olution de problèmes de mélange"
$s3 = "s"
$s4 = " grâce à une interface conviviale et facile d’utilisation"

$b1 = $p4Start + $s1.Length
$b2 = $b1 + $s2.Length
$b3 = $b2 + $s3.Length

$m1 = $d.Bookmarks.Add("ZZB1", $d.Range($b1, $b1))
$m2 = $d.Bookmarks.Add("ZZB2", $d.Range($b2, $b2))
$mg = $d.Bookmarks.Add("_GoBack", $d.Range($b3, $b3))

$r1 = $d.Range($p4Start, $b1)
$r1.Text = $r1.Text + "X"
$d.Range($b1, $b1 + 1).Delete()

Write-Output ("Para4: [" + $d.Paragraphs(4).Range.Text + "]")

$d.Bookmarks("ZZB1").Delete()
$d.Bookmarks("ZZB2").Delete()

Write-Output ("Para4 final: [" + $d.Paragraphs(4).Range.Text + "]")

# =====================================================================
# PART 4: insert a new empty paragraph (ind left=360) after paragraph 4
# =====================================================================

$p4 = $d.Paragraphs(4)
$null = $p4.Range.InsertParagraphAfter()

$p5 = $d.Paragraphs(5)
$p5.Range.ParagraphFormat.Style = "Normal"
$p5.Range.ParagraphFormat.LeftIndent = 18

$tmp = $d.Range($p5.Range.Start, $p5.Range.Start)
$tmp.InsertAfter("ZZZ")
$d.Range($p5.Range.Start, $p5.Range.Start + 3).Delete()

Write-Output ("Para5 (new): [" + $d.Paragraphs(5).Range.Text + "]")
Write-Output ("Para6: [" + $d.Paragraphs(6).Range.Text + "]")

Write-Output "DONE"
